# Weekly price update: a new week's records are prepended to the data
# block (rows 133-152), shifting the existing records down by two rows
# (to 135-154). The two newest records become the new rows 133-134.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the data block; everything from
# row 133 down (through the old last row 152) shifts down to 135-154.
$ws.Range("A133:R134").Insert()

# New row 133
$ws.Range("A133").Value = 4
$ws.Range("B133").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C133").Value = "Los Lagos"
$ws.Range("D133").Value = 44474
$ws.Range("E133").Value = 10
$ws.Range("F133").Value = 100112043
$ws.Range("G133").Value = "Pepino ensalada"
$ws.Range("H133").Value = "Sin especificar"
$ws.Range("I133").Value = "Primera"
$ws.Range("J133").Value = 200
$ws.Range("K133").Value = 23000
$ws.Range("L133").Value = 23000
$ws.Range("M133").Value = 23000
$ws.Range("N133").Value = "$/caja 60 unidades"
$ws.Range("O133").Value = "Región de Arica y Parinacota"
$ws.Range("P133").Value = 383
$ws.Range("Q133").Value = 60
$ws.Range("R133").Value = "Hortaliza"

# New row 134
$ws.Range("A134").Value = 4
$ws.Range("B134").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C134").Value = "Los Lagos"
$ws.Range("D134").Value = 44474
$ws.Range("E134").Value = 10
$ws.Range("F134").Value = 100112043
$ws.Range("G134").Value = "Pepino ensalada"
$ws.Range("H134").Value = "Sin especificar"
$ws.Range("I134").Value = "Segunda"
$ws.Range("J134").Value = 200
$ws.Range("K134").Value = 20000
$ws.Range("L134").Value = 20000
$ws.Range("M134").Value = 20000
$ws.Range("N134").Value = "$/caja 100 unidades"
$ws.Range("O134").Value = "Región de Arica y Parinacota"
$ws.Range("P134").Value = 200
$ws.Range("Q134").Value = 100
$ws.Range("R134").Value = "Hortaliza"
